$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; this shifts the existing B:J data right to C:K
# and widens the used range to A1:K26 automatically.
$ws.Columns("B").Insert()

# The insert doesn't carry over the custom column width, so restore it to
# match the other price columns (C:K, width 21).
$ws.Columns("B").ColumnWidth = $ws.Columns("C").ColumnWidth

# New (most-recent) timestamp header for the newly inserted column.
$ws.Range("B1").Value = "2025-12-21 10:29"

# Per-row price snapshot for the new 2025-12-21 10:29 column.
$newColB = @{
    2  = 929
    3  = 569
    4  = 299
    5  = 569
    6  = 499
    7  = 569
    8  = 929
    9  = 299
    10 = 299
    11 = 929
    12 = 569
    13 = 569
    14 = 499
    16 = 299
    17 = 929
    18 = 499
    19 = 1497
    20 = 929
    21 = 499
    22 = 299
    23 = 1299
    24 = 929
    25 = 929
    26 = 1299
}

foreach ($row in $newColB.Keys) {
    $ws.Cells.Item($row, 2).Value = $newColB[$row]
}

# Row 15's product had no price recorded at 2025-12-21 10:29 (out of stock).
$ws.Range("B15").Value = ""
